# The paragraph ends with runs:  ...amp; </del> </add></ab>
# (colors a9a9a9, a91111, 7f6000 respectively). The edit splits the last
# run ("</add></ab>", color 7f6000) into two pieces: a new "</add>" run
# that moves to sit *before* the "</del>" run, and the remainder
# "</ab>" which stays in place of the original run.

$d = $word.ActiveDocument

# 1. Grab the exact formatting (font, color, size...) of the run that
#    currently holds "</add></ab>" by capturing it as FormattedText -
#    this preserves the run's full rPr (rFonts ascii/eastAsia/hAnsi/cs,
#    color, sz, szCs, rtl) when we paste it elsewhere.
$srcRng = $d.Content
$null = $srcRng.Find.Execute("</add></ab>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$formatted = $d.Range($srcRng.Start, $srcRng.End).FormattedText

# 2. Locate the insertion point immediately before the "</del>" run that
#    precedes "</add></ab>" (use the longer needle so we land on the
#    right "</del>" - the document has several).
$target = $d.Content
$null = $target.Find.Execute("</del></add></ab>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertAt = $target.Start

# 3. Paste a full copy of "</add></ab>" (with its original formatting)
#    right before "</del>", then trim it down to just "</add>".
$insPoint = $d.Range($insertAt, $insertAt)
$insPoint.FormattedText = $formatted
$d.Range($insertAt + 6, $insertAt + 11).Delete()

# 4. Trim the original run (now after "</del>") from "</add></ab>" down
#    to "</ab>" - find it again since only one "</add></ab>" remains.
$orig = $d.Content
$null = $orig.Find.Execute("</add></ab>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Range($orig.Start, $orig.Start + 6).Delete()
